$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 94876
$ws.Range("E2").Value = 512
$ws.Range("F2").Value = 512
$ws.Range("G2").Value = -296
$ws.Range("H2").Value = -225
$ws.Range("I2").Value = -411
$ws.Range("J2").Value = 186
$ws.Range("K2").Value = 130947
$ws.Range("L2").Value = 95131
$ws.Range("M2").Value = 35816
$ws.Range("N2").Value = 32709
$ws.Range("O2").Value = 3107
$ws.Range("P2").Value = 3550
$ws.Range("Q2").Value = 5715
$ws.Range("R2").Value = -3345
$ws.Range("S2").Value = 352
$ws.Range("T2").Value = 3080
$ws.Range("U2").Value = 2636
$ws.Range("V2").Value = 39194
$ws.Range("W2").Value = 0.54
$ws.Range("X2").Value = -0.24
$ws.Range("Y2").Value = -1.34
$ws.Range("Z2").Value = -0.18
$ws.Range("AA2").Value = 265.61
$ws.Range("AB2").Value = 885.83
$ws.Range("AC2").Value = -649
$ws.Range("AD2").Value = -35.85
$ws.Range("AE2").Value = 47023
$ws.Range("AF2").Value = 0.49
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 71000000
$ws.Range("D3").Value = 105726
$ws.Range("E3").Value = 1221
$ws.Range("F3").Value = 1221
$ws.Range("G3").Value = 313
$ws.Range("H3").Value = 295
$ws.Range("I3").Value = 261
$ws.Range("J3").Value = 34
$ws.Range("K3").Value = 130055
$ws.Range("L3").Value = 96575
$ws.Range("M3").Value = 33480
$ws.Range("N3").Value = 32843
$ws.Range("O3").Value = 637
$ws.Range("P3").Value = 3550
$ws.Range("Q3").Value = -78
$ws.Range("R3").Value = 5079
$ws.Range("S3").Value = -2282
$ws.Range("T3").Value = 2271
$ws.Range("U3").Value = -2348
$ws.Range("V3").Value = 34621
$ws.Range("W3").Value = 1.16
$ws.Range("X3").Value = 0.28
$ws.Range("Y3").Value = 0.8
$ws.Range("Z3").Value = 0.23
$ws.Range("AA3").Value = 288.45
$ws.Range("AB3").Value = 886.75
$ws.Range("AC3").Value = 367
$ws.Range("AD3").Value = 53.8
$ws.Range("AE3").Value = 47215
$ws.Range("AF3").Value = 0.42
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 71000000
$ws.Range("D4").Value = 110356
$ws.Range("E4").Value = 1430
$ws.Range("F4").Value = 1430
$ws.Range("G4").Value = 213
$ws.Range("H4").Value = -204
$ws.Range("I4").Value = -258
$ws.Range("J4").Value = 54
$ws.Range("K4").Value = 133766
$ws.Range("L4").Value = 100234
$ws.Range("M4").Value = 33532
$ws.Range("N4").Value = 32788
$ws.Range("O4").Value = 744
$ws.Range("P4").Value = 3550
$ws.Range("Q4").Value = 812
$ws.Range("R4").Value = -2487
$ws.Range("S4").Value = 782
$ws.Range("T4").Value = 791
$ws.Range("U4").Value = 21
$ws.Range("V4").Value = 34335
$ws.Range("W4").Value = 1.3
$ws.Range("X4").Value = -0.19
$ws.Range("Y4").Value = -0.79
$ws.Range("Z4").Value = -0.16
$ws.Range("AA4").Value = 298.92
$ws.Range("AB4").Value = 877.48
$ws.Range("AC4").Value = -363
$ws.Range("AD4").Value = -72.98
$ws.Range("AE4").Value = 47136
$ws.Range("AF4").Value = 0.5600000000000001
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 71000000
$ws.Range("D5").Value = 116795
$ws.Range("E5").Value = 3187
$ws.Range("F5").Value = 3187
$ws.Range("G5").Value = -1607
$ws.Range("H5").Value = -1637
$ws.Range("I5").Value = -1684
$ws.Range("J5").Value = 47
$ws.Range("K5").Value = 136966
$ws.Range("L5").Value = 104569
$ws.Range("M5").Value = 32397
$ws.Range("N5").Value = 31673
$ws.Range("O5").Value = 724
$ws.Range("P5").Value = 3584
$ws.Range("Q5").Value = -2055
$ws.Range("R5").Value = 135
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 208
$ws.Range("U5").Value = -2263
$ws.Range("V5").Value = 39063
$ws.Range("W5").Value = 2.73
$ws.Range("X5").Value = -1.4
$ws.Range("Y5").Value = -5.22
$ws.Range("Z5").Value = -1.21
$ws.Range("AA5").Value = 322.77
$ws.Range("AB5").Value = 826.37
$ws.Range("AC5").Value = -2359
$ws.Range("AD5").Value = -12
$ws.Range("AE5").Value = 45095
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 1.06
$ws.Range("AI5").Value = -12.51
$ws.Range("AJ5").Value = 71675237
$ws.Range("D6").Value = 131394
$ws.Range("E6").Value = 10645
$ws.Range("F6").Value = 10645
$ws.Range("G6").Value = 8350
$ws.Range("H6").Value = 5874
$ws.Range("I6").Value = 5821
$ws.Range("K6").Value = 120333
$ws.Range("L6").Value = 84086
$ws.Range("M6").Value = 36247
$ws.Range("N6").Value = 35485
$ws.Range("P6").Value = 3972
$ws.Range("Q6").Value = 10329
$ws.Range("R6").Value = -2348
$ws.Range("S6").Value = -16556
$ws.Range("T6").Value = 137
$ws.Range("U6").Value = 10193
$ws.Range("V6").Value = 21030
$ws.Range("W6").Value = 8.1
$ws.Range("X6").Value = 4.47
$ws.Range("Y6").Value = 17.34
$ws.Range("Z6").Value = 4.57
$ws.Range("AA6").Value = 231.98
$ws.Range("AB6").Value = 830.5700000000001
$ws.Range("AC6").Value = 7713
$ws.Range("AD6").Value = 5.67
$ws.Range("AE6").Value = 45064
$ws.Range("AF6").Value = 0.97
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 2.29
$ws.Range("AI6").Value = 13.53
$ws.Range("AJ6").Value = 79435797
$ws.Range("D7").Value = 102908
$ws.Range("E7").Value = 7669
$ws.Range("G7").Value = 7885
$ws.Range("H7").Value = 5673
$ws.Range("I7").Value = 5600
$ws.Range("K7").Value = 127765
$ws.Range("L7").Value = 86230
$ws.Range("M7").Value = 41535
$ws.Range("N7").Value = 40606
$ws.Range("P7").Value = 3999
$ws.Range("Q7").Value = 5208
$ws.Range("R7").Value = -2800
$ws.Range("S7").Value = -922
$ws.Range("T7").Value = 1147
$ws.Range("U7").Value = 4174
$ws.Range("W7").Value = 7.45
$ws.Range("X7").Value = 5.51
$ws.Range("Y7").Value = 14.72
$ws.Range("Z7").Value = 4.57
$ws.Range("AA7").Value = 207.61
$ws.Range("AC7").Value = 7012
$ws.Range("AD7").Value = 4.32
$ws.Range("AE7").Value = 51137
$ws.Range("AF7").Value = 0.59
$ws.Range("AG7").Value = 1024
$ws.Range("AH7").Value = 3.38
$ws.Range("AI7").Value = 14.64
$ws.Range("D8").Value = 104247
$ws.Range("E8").Value = 7447
$ws.Range("G8").Value = 6916
$ws.Range("H8").Value = 5068
$ws.Range("I8").Value = 5014
$ws.Range("K8").Value = 130662
$ws.Range("L8").Value = 84801
$ws.Range("M8").Value = 45861
$ws.Range("N8").Value = 44961
$ws.Range("P8").Value = 4019
$ws.Range("Q8").Value = 6011
$ws.Range("R8").Value = -2130
$ws.Range("S8").Value = -3204
$ws.Range("T8").Value = 532
$ws.Range("U8").Value = 4518
$ws.Range("W8").Value = 7.14
$ws.Range("X8").Value = 4.86
$ws.Range("Y8").Value = 11.72
$ws.Range("Z8").Value = 3.92
$ws.Range("AA8").Value = 184.91
$ws.Range("AC8").Value = 6260
$ws.Range("AD8").Value = 4.47
$ws.Range("AE8").Value = 56622
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 1075
$ws.Range("AH8").Value = 3.84
$ws.Range("AI8").Value = 17.17
$ws.Range("D9").Value = 107346
$ws.Range("E9").Value = 7704
$ws.Range("G9").Value = 7352
$ws.Range("H9").Value = 5370
$ws.Range("I9").Value = 5425
$ws.Range("K9").Value = 135633
$ws.Range("L9").Value = 85252
$ws.Range("M9").Value = 50381
$ws.Range("N9").Value = 49636
$ws.Range("P9").Value = 4019
$ws.Range("Q9").Value = 6556
$ws.Range("R9").Value = -2527
$ws.Range("S9").Value = -3283
$ws.Range("T9").Value = 536
$ws.Range("U9").Value = 4644
$ws.Range("W9").Value = 7.18
$ws.Range("X9").Value = 5
$ws.Range("Y9").Value = 11.47
$ws.Range("Z9").Value = 4.03
$ws.Range("AA9").Value = 169.22
$ws.Range("AC9").Value = 6773
$ws.Range("AD9").Value = 4.13
$ws.Range("AE9").Value = 62510
$ws.Range("AF9").Value = 0.45
$ws.Range("AG9").Value = 1088
$ws.Range("AH9").Value = 3.88
$ws.Range("AI9").Value = 16.06
